# Add team Win/Loss/Tie record columns (AD:AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers Wins / Losses / Ties ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold font, thin border,
# centered/top aligned) by copying the existing header style onto the new
# header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2-38): team record values ---
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
